# Actualización automática 2025-06-16 13:01:14
# Adds a "GRANITO" product column (inserted before GRIFERIAS) and three
# trailing product columns (NO RESURTIBLES, PANELES PVC, PANELES PU) to the
# "VENTAS POR GRUPO" sheet, extending the table from A1:N5 to A1:R5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new "GRANITO" column before column F (GRIFERIAS) ---------
# Inserting a full column shifts the existing F:N columns to G:O and
# automatically copies cell formatting/styles from the neighbouring column,
# so header (s=1), data (s=2) and summary-row (s=3) styles carry over.
$ws.Columns("F:F").Insert()

$ws.Range("F1").Value = "GRANITO"
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = "0 de 3"

# Restore the intended column width for the new column (13 characters).
# ColumnWidth and the stored XML "width" differ by a constant offset
# (~0.83) caused by Excel's character-to-pixel padding, so compensate here.
$ws.Columns("F:F").ColumnWidth = 13 - 0.83

# --- Append three new trailing columns: P, Q, R ---------------------------
# Copy formatting (header/data/summary styles) from the last existing
# column (O, formerly N = "SAL SOLUBLE") across the three new columns.
$ws.Range("O1:O5").Copy()
$ws.Range("P1:R5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("P1").Value = "NO RESURTIBLES"
$ws.Range("Q1").Value = "PANELES PVC"
$ws.Range("R1").Value = "PANELES PU"

$ws.Range("P2:R4").Value = 0

$ws.Range("P5").Value = "0 de 3"
$ws.Range("Q5").Value = "0 de 3"
$ws.Range("R5").Value = "0 de 3"

$ws.Columns("P:P").ColumnWidth = 20 - 0.83
$ws.Columns("Q:Q").ColumnWidth = 17 - 0.83
$ws.Columns("R:R").ColumnWidth = 16 - 0.83
